$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")
[void]$ws.Select()

# Insert two new rows right above the current row 4 ("dataset.commit.id"),
# pushing all following rows down by two.
[void]$ws.Rows.Item(4).Insert()
[void]$ws.Rows.Item(4).Insert()

$tableFormula = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"

$lineFormula = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = $tableFormula
$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = $lineFormula

$newRange = $ws.Range("A4:B5")
$newRange.VerticalAlignment = -4108
$newRange.WrapText = $true

$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(5).RowHeight = 120

[void]$ws.Range("B10").Select()
